$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(17).Insert()
$ws.Rows("17:17").RowHeight = 13.2
